$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: "sheet1" -> "sheet1 -1"
$ws.Name = "sheet1 -1"

# Widen column B slightly (53.25 -> ~56.58 chars)
$ws.Columns.Item(2).ColumnWidth = 55.79

# Append three new keyword/appID rows
$ws.Range("A10").Value = "passive"
$ws.Range("B10").Value = "passive.income.nadi.myfirstdrawermenuproject"
$ws.Range("A11").Value = "income"
$ws.Range("B11").Value = "passive.income.nadi.myfirstdrawermenuproject"
$ws.Range("A12").Value = "stretchy"
$ws.Range("B12").Value = "com.singleton.strechy"

# Carry over the same cell style (wrap text, etc.) used by the rest of the table
$ws.Range("A9:B9").Copy()
$ws.Range("A10:B12").PasteSpecial(-4122)

# Row 10 wraps onto two lines (long appID), same as rows 3 & 7
$ws.Rows.Item(10).RowHeight = 24

# Move the active selection down to the new last row
$ws.Range("A12").Select() | Out-Null
